$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '72.170.84'
$ws.Cells.Item(2, 5).Value = '  +1.69%  '
$ws.Cells.Item(3, 4).Value = '2.668.19'
$ws.Cells.Item(3, 5).Value = '  +2.12%  '
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$ws.Cells.Item(5, 4).Value = '602.18'
$ws.Cells.Item(5, 5).Value = '  -0.39%  '
$ws.Cells.Item(6, 4).Value = '178.63'
$ws.Cells.Item(6, 5).Value = '  -0.81%  '
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
$ws.Cells.Item(8, 5).Value = '  -0.67%  '
$ws.Cells.Item(9, 4).Value = '2.666.01'
$ws.Cells.Item(9, 5).Value = '  +2.08%  '
$ws.Cells.Item(10, 5).Value = '  +4.98%  '
$ws.Cells.Item(11, 5).Value = '  +2.08%  '
$ws.Cells.Item(12, 5).Value = '  +2.77%  '
$ws.Cells.Item(13, 4).Value = '5.03'
$ws.Cells.Item(13, 5).Value = '  +0.26%  '
$ws.Cells.Item(14, 4).Value = '3.153.06'
$ws.Cells.Item(14, 5).Value = '  +1.77%  '
$ws.Cells.Item(15, 5).Value = '  +2.95%  '
$ws.Cells.Item(16, 4).Value = '72.047.65'
$ws.Cells.Item(16, 5).Value = '  +1.50%  '
$ws.Cells.Item(17, 4).Value = '26.59'
$ws.Cells.Item(17, 5).Value = '  +0.16%  '
$ws.Cells.Item(18, 4).Value = '2.670.05'
$ws.Cells.Item(18, 5).Value = '  +2.15%  '
$ws.Cells.Item(19, 4).Value = '11.97'
$ws.Cells.Item(19, 5).Value = '  +4.28%  '
$ws.Cells.Item(20, 4).Value = '8.01'
$ws.Cells.Item(20, 5).Value = '  +3.51%  '
$ws.Cells.Item(21, 4).Value = '380.72'
$ws.Cells.Item(21, 5).Value = '  +0.01%  '
$ws.Cells.Item(22, 5).Value = '  +1.46%  '
$ws.Cells.Item(23, 5).Value = '  +11.31%  '
$ws.Cells.Item(24, 4).Value = '72.44'
$ws.Cells.Item(24, 5).Value = '  +0.31%  '
$ws.Cells.Item(25, 5).Value = '  -0.12%  '
$ws.Cells.Item(26, 5).Value = '  -1.49%  '
$ws.Cells.Item(27, 5).Value = '  +3.93%  '
$ws.Cells.Item(28, 4).Value = '2.809.60'
$ws.Cells.Item(28, 5).Value = '  +3.59%  '
$ws.Cells.Item(29, 4).Value = '0.998'
$ws.Cells.Item(29, 5).Value = '  -0.03%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0947'
$ws.Cells.Item(30, 5).Value = '  -0.13%  '
$ws.Cells.Item(31, 4).Value = '8.16'
$ws.Cells.Item(31, 5).Value = '  +1.89%  '
$ws.Cells.Item(32, 4).Value = '520.86'
$ws.Cells.Item(32, 5).Value = '  -1.11%  '
$ws.Cells.Item(33, 5).Value = '  -0.39%  '
$ws.Cells.Item(34, 5).Value = '  -0.49%  '
$ws.Cells.Item(35, 5).Value = '  -0.04%  '
$ws.Cells.Item(36, 4).Value = '164.68'
$ws.Cells.Item(36, 5).Value = '  +0.43%  '
$ws.Cells.Item(37, 4).Value = '19.57'
$ws.Cells.Item(37, 5).Value = '  +2.25%  '
$ws.Cells.Item(38, 5).Value = '  +0.86%  '
$ws.Cells.Item(39, 5).Value = '  +1.62%  '
$ws.Cells.Item(40, 5).Value = '  -6.67%  '
$ws.Cells.Item(41, 5).Value = '  -1.84%  '
$ws.Cells.Item(43, 5).Value = '  +0.64%  '
$ws.Cells.Item(44, 5).Value = '  -0.75%  '
$ws.Cells.Item(45, 5).Value = '  +1.75%  '
$ws.Cells.Item(46, 4).Value = '39.29'
$ws.Cells.Item(46, 5).Value = '  -1.82%  '
$ws.Cells.Item(47, 4).Value = '152.78'
$ws.Cells.Item(47, 5).Value = '  -0.85%  '
$ws.Cells.Item(48, 5).Value = '  +3.36%  '
$ws.Cells.Item(49, 5).Value = '  +3.66%  '
$ws.Cells.Item(50, 5).Value = '  +3.05%  '
$ws.Cells.Item(51, 4).Value = '0.0765'
$ws.Cells.Item(51, 5).Value = '  +1.68%  '
